$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 999.8
$ws.Range("I32").Value = 999.8
$ws.Range("K32").Value = 999.8
$ws.Range("M32").Value = -673.8

# Row 138
$ws.Range("H138").Value = 2887.5
$ws.Range("I138").Value = 2887
$ws.Range("K138").Value = 8661
$ws.Range("M138").Value = -3521

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1522.5
$ws.Range("I2").Value = 1522.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1522.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1409.5
$ws.Range("N2").ClearContents()

# Row 88
$ws.Range("H88").Value = 550
$ws.Range("I88").Value = 550
$ws.Range("K88").Value = 550
$ws.Range("M88").Value = -144

# Row 91
$ws.Range("H91").Value = 550
$ws.Range("I91").Value = 550
$ws.Range("K91").Value = 550
$ws.Range("M91").Value = 854

# Row 97
$ws.Range("H97").Value = 412.75
$ws.Range("J97").Value = 716.3333
$ws.Range("L97").Value = 716.3333
$ws.Range("N97").Value = -1708.3333

# Row 116
$ws.Range("H116").Value = 1522.5
$ws.Range("I116").Value = 1522.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1522.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 771.5
$ws.Range("N116").ClearContents()

# Row 122
$ws.Range("H122").Value = 1703.2
$ws.Range("I122").Value = 1703.2
$ws.Range("K122").Value = 5109.6
$ws.Range("M122").Value = -2659.6

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1522.5
$ws.Range("I3").Value = 1522.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1522.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1408.5
$ws.Range("N3").ClearContents()

# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# Row 86
$ws.Range("H86").Value = 3280.524
$ws.Range("I86").Value = 3031.3125
$ws.Range("J86").Value = 4078
$ws.Range("K86").Value = 3031.3125
$ws.Range("L86").Value = 4078
$ws.Range("M86").Value = -1908.3125
$ws.Range("N86").Value = -6324

# Row 89
$ws.Range("H89").Value = 3280.524
$ws.Range("I89").Value = 3031.3125
$ws.Range("J89").Value = 4078
$ws.Range("K89").Value = 15156.5625
$ws.Range("L89").Value = 20390
$ws.Range("M89").Value = -9540.5625
$ws.Range("N89").Value = -31622

# Row 94
$ws.Range("H94").Value = 1403.5714
$ws.Range("I94").Value = 1470.8334
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 1470.8334
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -1019.8334
$ws.Range("N94").Value = -1902

# Row 96
$ws.Range("H96").Value = 19500
$ws.Range("I96").Value = 19500
$ws.Range("K96").Value = 19500
$ws.Range("M96").Value = -16754

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 313.5
$ws.Range("I10").Value = 265.875
$ws.Range("K10").Value = 265.875
$ws.Range("M10").Value = -126.875

# Row 13
$ws.Range("H13").Value = 5000
$ws.Range("I13").Value = 5000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -4861
$ws.Range("N13").ClearContents()

# Row 132
$ws.Range("H132").Value = 3645.25
$ws.Range("I132").Value = 3889.2222
$ws.Range("K132").Value = 11667.6666
$ws.Range("M132").Value = -9137.6666

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 21.541666
$ws.Range("I2").Value = 18.9375
$ws.Range("J2").Value = 26.75
$ws.Range("K2").Value = 113.625
$ws.Range("L2").Value = 160.5
$ws.Range("M2").Value = -0.625
$ws.Range("N2").Value = -386.5

# Row 69
$ws.Range("H69").Value = 2571.3333
$ws.Range("J69").Value = 2571.3333
$ws.Range("L69").Value = 7713.999899999999
$ws.Range("N69").Value = -9335.999899999999

# Row 72
$ws.Range("H72").Value = 2571.3333
$ws.Range("J72").Value = 2571.3333
$ws.Range("L72").Value = 23141.9997
$ws.Range("N72").Value = -31253.9997

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 892.375
$ws.Range("I97").Value = 805.5714
$ws.Range("K97").Value = 805.5714
$ws.Range("M97").Value = -309.5714

$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 20
$ws.Range("H20").Value = 16667
$ws.Range("J20").Value = 16667
$ws.Range("L20").Value = 16667
$ws.Range("N20").Value = -17119

# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

# Row 22
$ws.Range("H22").Value = 595.9167
$ws.Range("J22").Value = 642.6667
$ws.Range("L22").Value = 642.6667
$ws.Range("N22").Value = -1232.6667

# Row 27
$ws.Range("H27").Value = 595.9167
$ws.Range("J27").Value = 642.6667
$ws.Range("L27").Value = 642.6667
$ws.Range("N27").Value = -856.6667

# Row 68
$ws.Range("H68").Value = 1311.875
$ws.Range("I68").Value = 915.8333
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 915.8333
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -166.8333
$ws.Range("N68").Value = -3998

# Row 71
$ws.Range("H71").Value = 1311.875
$ws.Range("I71").Value = 915.8333
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 4579.1665
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -835.1665000000003
$ws.Range("N71").Value = -19988

# Row 93
$ws.Range("H93").Value = 1400
$ws.Range("I93").Value = 1400
$ws.Range("K93").Value = 1400
$ws.Range("M93").Value = -152

# Row 136
$ws.Range("H136").Value = 1797.3334
$ws.Range("I136").Value = 1797.3334
$ws.Range("K136").Value = 5392.0002
$ws.Range("M136").Value = -2842.0002

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 14500
$ws.Range("I2").Value = 20000
$ws.Range("K2").Value = 20000
$ws.Range("M2").Value = -19888

# Row 6
$ws.Range("H6").Value = 851.25
$ws.Range("I6").Value = 502.5
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 502.5
$ws.Range("L6").Value = 1200
$ws.Range("M6").Value = -387.5
$ws.Range("N6").Value = -1430

# Row 100
$ws.Range("H100").Value = 346
$ws.Range("I100").Value = 346
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 692
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -151
$ws.Range("N100").ClearContents()
